# Updates cryptos list values (Price / Volume(1h)) to match the
# "Updated cryptos list" data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.257.35'
$ws.Range("E2").Value = '  +1.81%  '
$ws.Range("D3").Value = '2.460.12'
$ws.Range("E3").Value = '  +2.50%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''565.58'
$ws.Range("E5").Value = '  +1.05%  '
$ws.Range("D6").Value = '''142.97'
$ws.Range("E6").Value = '  +3.86%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '''0.584'
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("D9").Value = '2.459.60'
$ws.Range("E9").Value = '  +2.55%  '
$ws.Range("D10").Value = '''0.105'
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("D11").Value = '''5.70'
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("E12").Value = '  +1.60%  '
$ws.Range("D13").Value = '''0.356'
$ws.Range("E13").Value = '  +1.82%  '
$ws.Range("D14").Value = '''27.21'
$ws.Range("E14").Value = '  +5.29%  '
$ws.Range("D15").Value = '2.898.62'
$ws.Range("E15").Value = '  +2.47%  '
$ws.Range("D16").Value = '63.053.62'
$ws.Range("E16").Value = '  +1.56%  '
$ws.Range("D17").Value = '''0.0000141'
$ws.Range("E17").Value = '  +1.63%  '
$ws.Range("D18").Value = '2.456.04'
$ws.Range("E18").Value = '  +0.86%  '
$ws.Range("D19").Value = '''11.27'
$ws.Range("E19").Value = '  +2.00%  '
$ws.Range("D20").Value = '''340.72'
$ws.Range("E20").Value = '  -1.22%  '
$ws.Range("D21").Value = '''4.29'
$ws.Range("E21").Value = '  +1.41%  '
$ws.Range("E22").Value = '  -2.02%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '''65.58'
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("D25").Value = '''0.171'
$ws.Range("E25").Value = '  -1.75%  '
$ws.Range("D26").Value = '''0.997'
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").Value = '''1.51'
$ws.Range("E27").Value = '  +0.91%  '
$ws.Range("E28").Value = '  +4.82%  '
$ws.Range("D29").Value = '''8.08'
$ws.Range("E29").Value = '  -3.03%  '
$ws.Range("D32").Value = '0.0₃0792'
$ws.Range("E32").Value = '  +3.54%  '
$ws.Range("D33").Value = '''176.24'
$ws.Range("E33").Value = '  +2.87%  '
$ws.Range("D34").Value = '''1.52'
$ws.Range("E34").Value = '  +7.91%  '
$ws.Range("D35").Value = '''389.84'
$ws.Range("E35").Value = '  +9.00%  '
$ws.Range("D36").Value = '''0.398'
$ws.Range("E36").Value = '  +1.19%  '
$ws.Range("D37").Value = '''18.77'
$ws.Range("E37").Value = '  +1.18%  '
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("D39").Value = '''4.33'
$ws.Range("E39").Value = '  -4.51%  '
$ws.Range("D40").Value = '''1.74'
$ws.Range("E40").Value = '  +4.32%  '
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("D42").Value = '''40.03'
$ws.Range("E42").Value = '  +2.49%  '
$ws.Range("D43").Value = '''149.71'
$ws.Range("E43").Value = '  +4.13%  '
$ws.Range("D44").Value = '''3.70'
$ws.Range("E44").Value = '  +0.70%  '
$ws.Range("D45").Value = '''20.63'
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("D46").Value = '''0.598'
$ws.Range("E46").Value = '  +2.57%  '
$ws.Range("D47").Value = '''0.0962'
$ws.Range("E47").Value = '  -0.24%  '
$ws.Range("D48").Value = '''0.0516'
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("D49").Value = '''0.0229'
$ws.Range("E49").Value = '  +3.50%  '
$ws.Range("D50").Value = '0.0₆0230'
$ws.Range("E50").Value = '  +7.05%  '
$ws.Range("D51").Value = '''17.95'
$ws.Range("E51").Value = '  +0.64%  '

# Rows 30/31 swapped coins (Aptos <-> PancakeSwap) along with new values
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '''1.85'
$ws.Range("E30").Value = '  +2.66%  '
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").Value = '''6.76'
$ws.Range("E31").Value = '  +6.36%  '
